# Insert a new price-report row for "Femacal de La Calera - Poroto verde"
# at row 148, pushing the existing rows 148:213 down to 149:214.
# (Sheet dimension grows from A1:R213 to A1:R214.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148 and below down by one row.
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with the new data point.
$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44466
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112031
$ws.Range("G148").Value = "Poroto verde"
$ws.Range("H148").Value = "Magnum"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 38
$ws.Range("K148").Value = 35000
$ws.Range("L148").Value = 35000
$ws.Range("M148").Value = 35000
$ws.Range("N148").Value = "$/malla 25 kilos"
$ws.Range("O148").Value = "Región de Arica y Parinacota"
$ws.Range("P148").Value = 1400
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = "Hortaliza"
